# Actualización automática 2025-06-25 17:25:09
$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M5").Value = 3542.21
$ws1.Range("H10").Value = 2376
$ws1.Range("I10").Value = 486
$ws1.Range("L29").Value = 2315.36
$ws1.Range("E31").Value = 381.84
$ws1.Range("E54").Value = "3 de 52"
$ws1.Range("H54").Value = "3 de 52"
$ws1.Range("I54").Value = "3 de 52"

# --- Sheet: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F5").Value = 4283.91
$ws2.Range("F10").Value = 2862
$ws2.Range("F29").Value = 6039.12
$ws2.Range("F31").Value = 381.84
$ws2.Range("F54").Value = 65381.22

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws3.Range("D4").Value = 1101.5
$ws3.Range("E4").Value = -98.5
$ws3.Range("F4").Value = 1.098205383848455

$ws3.Range("D7").Value = 5118.3
$ws3.Range("E7").Value = -3718.3
$ws3.Range("F7").Value = 3.655928571428571

$ws3.Range("D8").Value = 1019.63
$ws3.Range("E8").Value = -19.63
$ws3.Range("F8").Value = 1.01963

$ws3.Range("D15").Value = 19919.42
$ws3.Range("E15").Value = -6419.419999999998
$ws3.Range("F15").Value = 1.475512592592592

$ws3.Range("D16").Value = 14584.97
$ws3.Range("E16").Value = 18156.48
$ws3.Range("F16").Value = 0.4454588907943905

$ws3.Range("D19").Value = 67796.14999999999
$ws3.Range("E19").Value = 26651.29064517915
$ws3.Range("F19").Value = 0.7178188158078003
